$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update IPC PO (C), DELTA (D) and DELTA^2 (E) results for the partial model's
# sliding-window predictions (rows 2-51), plus a handful of IPC RO (B) cells whose
# stored double differs in its last bit, and the TOTAL/MSE summary cells
# (C52 = SUM(DELTA), E52 = SUM(DELTA^2), E53 = MSE = AVERAGE(DELTA^2)).

$ws.Cells.Item(2, 3).Value = 29.30314445495605
$ws.Cells.Item(2, 4).Value = 0.3231444549560507
$ws.Cells.Item(2, 5).Value = 0.1044223387688431
$ws.Cells.Item(3, 2).Value = 29.15000000000001
$ws.Cells.Item(3, 3).Value = 29.25433921813965
$ws.Cells.Item(3, 4).Value = 0.1043392181396428
$ws.Cells.Item(3, 5).Value = 0.01088667244199196
$ws.Cells.Item(4, 2).Value = 29.34999999999999
$ws.Cells.Item(4, 3).Value = 29.70592308044434
$ws.Cells.Item(4, 4).Value = 0.3559230804443416
$ws.Cells.Item(4, 5).Value = 0.1266812391929893
$ws.Cells.Item(5, 3).Value = 29.52533531188965
$ws.Cells.Item(5, 4).Value = 0.1553353118896439
$ws.Cells.Item(5, 5).Value = 0.02412905911985294
$ws.Cells.Item(6, 2).Value = 29.53999999999999
$ws.Cells.Item(6, 3).Value = 28.87110710144043
$ws.Cells.Item(6, 4).Value = -0.6688928985595624
$ws.Cells.Item(6, 5).Value = 0.447417709743413
$ws.Cells.Item(7, 3).Value = 29.35497856140137
$ws.Cells.Item(7, 4).Value = -0.19502143859863
$ws.Cells.Item(7, 5).Value = 0.0380333615130792
$ws.Cells.Item(8, 3).Value = 29.55397605895996
$ws.Cells.Item(8, 4).Value = -0.1960239410400391
$ws.Cells.Item(8, 5).Value = 0.03842538546086871
$ws.Cells.Item(9, 3).Value = 30.04802513122559
$ws.Cells.Item(9, 4).Value = 0.2080251312255825
$ws.Cells.Item(9, 5).Value = 0.04327445522142083
$ws.Cells.Item(10, 3).Value = 29.96916007995605
$ws.Cells.Item(10, 4).Value = 0.1591600799560524
$ws.Cells.Item(10, 5).Value = 0.025331931051617
$ws.Cells.Item(11, 3).Value = 29.92043113708496
$ws.Cells.Item(11, 4).Value = 0.0004311370849592322
$ws.Cells.Item(11, 5).Value = 0.0000001858791860271442
$ws.Cells.Item(12, 3).Value = 29.84786033630371
$ws.Cells.Item(12, 4).Value = -0.132139663696293
$ws.Cells.Item(12, 5).Value = 0.01746089072176943
$ws.Cells.Item(13, 2).Value = 30.03999999999999
$ws.Cells.Item(13, 3).Value = 30.05769920349121
$ws.Cells.Item(13, 4).Value = 0.0176992034912189
$ws.Cells.Item(13, 5).Value = 0.0003132618042235752
$ws.Cells.Item(14, 2).Value = 30.21000000000001
$ws.Cells.Item(14, 3).Value = 30.05219078063965
$ws.Cells.Item(14, 4).Value = -0.1578092193603595
$ws.Cells.Item(14, 5).Value = 0.02490374971512607
$ws.Cells.Item(15, 3).Value = 30.22921562194824
$ws.Cells.Item(15, 4).Value = 0.009215621948243324
$ws.Cells.Item(15, 5).Value = 0.00008492768789294409
$ws.Cells.Item(16, 3).Value = 30.29749870300293
$ws.Cells.Item(16, 4).Value = -0.08250129699706577
$ws.Cells.Item(16, 5).Value = 0.006806464006198052
$ws.Cells.Item(17, 3).Value = 30.6120433807373
$ws.Cells.Item(17, 4).Value = 0.172043380737307
$ws.Cells.Item(17, 5).Value = 0.02959892485552196
$ws.Cells.Item(18, 3).Value = 30.40457725524902
$ws.Cells.Item(18, 4).Value = -0.07542274475098054
$ws.Cells.Item(18, 5).Value = 0.005688590425771562
$ws.Cells.Item(19, 3).Value = 30.41995811462402
$ws.Cells.Item(19, 4).Value = -0.2700418853759743
$ws.Cells.Item(19, 5).Value = 0.07292261985741083
$ws.Cells.Item(20, 3).Value = 30.54005241394043
$ws.Cells.Item(20, 4).Value = -0.2099475860595703
$ws.Cells.Item(20, 5).Value = 0.04407798889224068
$ws.Cells.Item(21, 3).Value = 30.62829399108887
$ws.Cells.Item(21, 4).Value = -0.3117060089111305
$ws.Cells.Item(21, 5).Value = 0.0971606359913058
$ws.Cells.Item(22, 3).Value = 30.73287773132324
$ws.Cells.Item(22, 4).Value = -0.2171222686767607
$ws.Cells.Item(22, 5).Value = 0.04714207955534344
$ws.Cells.Item(23, 3).Value = 31.12346458435059
$ws.Cells.Item(23, 4).Value = 0.1034645843505899
$ws.Cells.Item(23, 5).Value = 0.01070492021484034
$ws.Cells.Item(24, 3).Value = 31.30785179138184
$ws.Cells.Item(24, 4).Value = 0.1878517913818314
$ws.Cells.Item(24, 5).Value = 0.0352882955253631
$ws.Cells.Item(25, 3).Value = 31.35338401794434
$ws.Cells.Item(25, 4).Value = 0.0733840179443348
$ws.Cells.Item(25, 5).Value = 0.005385214089654452
$ws.Cells.Item(26, 3).Value = 31.22195243835449
$ws.Cells.Item(26, 4).Value = -0.1580475616455033
$ws.Cells.Item(26, 5).Value = 0.02497903174208916
$ws.Cells.Item(27, 3).Value = 31.44917106628418
$ws.Cells.Item(27, 4).Value = -0.1308289337158186
$ws.Cells.Item(27, 5).Value = 0.01711620989721806
$ws.Cells.Item(28, 2).Value = 31.65000000000001
$ws.Cells.Item(28, 3).Value = 31.89203453063965
$ws.Cells.Item(28, 4).Value = 0.2420345306396428
$ws.Cells.Item(28, 5).Value = 0.05858071402195217
$ws.Cells.Item(29, 3).Value = 32.54359817504883
$ws.Cells.Item(29, 4).Value = 0.6635981750488327
$ws.Cells.Item(29, 5).Value = 0.4403625379281412
$ws.Cells.Item(30, 3).Value = 32.46515274047852
$ws.Cells.Item(30, 4).Value = 0.1851527404785145
$ws.Cells.Item(30, 5).Value = 0.03428153730670414
$ws.Cells.Item(31, 3).Value = 32.5944709777832
$ws.Cells.Item(31, 4).Value = 0.1444709777832003
$ws.Cells.Item(31, 5).Value = 0.02087186342163395
$ws.Cells.Item(32, 2).Value = 32.84999999999999
$ws.Cells.Item(32, 3).Value = 32.88671493530273
$ws.Cells.Item(32, 4).Value = 0.03671493530274006
$ws.Cells.Item(32, 5).Value = 0.001347986474284388
$ws.Cells.Item(33, 2).Value = 32.90000000000001
$ws.Cells.Item(33, 3).Value = 33.03625106811523
$ws.Cells.Item(33, 4).Value = 0.1362510681152287
$ws.Cells.Item(33, 5).Value = 0.01856435356254069
$ws.Cells.Item(34, 2).Value = 33.09999999999999
$ws.Cells.Item(34, 3).Value = 32.95847702026367
$ws.Cells.Item(34, 4).Value = -0.1415229797363224
$ws.Cells.Item(34, 5).Value = 0.02002875379344753
$ws.Cells.Item(35, 2).Value = 33.40000000000001
$ws.Cells.Item(35, 3).Value = 33.7332878112793
$ws.Cells.Item(35, 4).Value = 0.3332878112792912
$ws.Cells.Item(35, 5).Value = 0.1110807651473404
$ws.Cells.Item(36, 3).Value = 33.6826286315918
$ws.Cells.Item(36, 4).Value = -0.01737136840820597
$ws.Cells.Item(36, 5).Value = 0.0003017644403736163
$ws.Cells.Item(37, 2).Value = 34.09999999999999
$ws.Cells.Item(37, 3).Value = 33.8400993347168
$ws.Cells.Item(37, 4).Value = -0.2599006652831974
$ws.Cells.Item(37, 5).Value = 0.06754835581464863
$ws.Cells.Item(38, 2).Value = 34.40000000000001
$ws.Cells.Item(38, 3).Value = 34.36728286743164
$ws.Cells.Item(38, 4).Value = -0.03271713256836506
$ws.Cells.Item(38, 5).Value = 0.001070410763495974
$ws.Cells.Item(39, 2).Value = 34.90000000000001
$ws.Cells.Item(39, 3).Value = 34.92761611938477
$ws.Cells.Item(39, 4).Value = 0.02761611938475994
$ws.Cells.Item(39, 5).Value = 0.0007626500498733137
$ws.Cells.Item(40, 3).Value = 35.70381546020508
$ws.Cells.Item(40, 4).Value = 0.403815460205081
$ws.Cells.Item(40, 5).Value = 0.1630669259006413
$ws.Cells.Item(41, 3).Value = 35.95959854125977
$ws.Cells.Item(41, 4).Value = 0.2595985412597628
$ws.Cells.Item(41, 5).Value = 0.06739140262419677
$ws.Cells.Item(42, 3).Value = 36.04191970825195
$ws.Cells.Item(42, 4).Value = -0.258080291748044
$ws.Cells.Item(42, 5).Value = 0.06660543698875553
$ws.Cells.Item(43, 3).Value = 36.59153366088867
$ws.Cells.Item(43, 4).Value = -0.2084663391113253
$ws.Cells.Item(43, 5).Value = 0.04345821454247807
$ws.Cells.Item(44, 3).Value = 37.04032516479492
$ws.Cells.Item(44, 4).Value = -0.2596748352050753
$ws.Cells.Item(44, 5).Value = 0.067431020038783
$ws.Cells.Item(45, 2).Value = 37.90000000000001
$ws.Cells.Item(45, 3).Value = 37.8508186340332
$ws.Cells.Item(45, 4).Value = -0.04918136596680256
$ws.Cells.Item(45, 5).Value = 0.002418806758360565
$ws.Cells.Item(46, 3).Value = 38.36010360717773
$ws.Cells.Item(46, 4).Value = -0.1398963928222656
$ws.Cells.Item(46, 5).Value = 0.01957100072468165
$ws.Cells.Item(47, 2).Value = 38.90000000000001
$ws.Cells.Item(47, 3).Value = 38.94926071166992
$ws.Cells.Item(47, 4).Value = 0.04926071166991619
$ws.Cells.Item(47, 5).Value = 0.002426617714226617
$ws.Cells.Item(48, 2).Value = 39.40000000000001
$ws.Cells.Item(48, 3).Value = 39.48398208618164
$ws.Cells.Item(48, 4).Value = 0.08398208618163494
$ws.Cells.Item(48, 5).Value = 0.007052990799419559
$ws.Cells.Item(49, 2).Value = 39.90000000000001
$ws.Cells.Item(49, 3).Value = 39.5393180847168
$ws.Cells.Item(49, 4).Value = -0.3606819152832088
$ws.Cells.Item(49, 5).Value = 0.1300914440123638
$ws.Cells.Item(50, 2).Value = 40.09999999999999
$ws.Cells.Item(50, 3).Value = 40.14670181274414
$ws.Cells.Item(50, 4).Value = 0.04670181274414631
$ws.Cells.Item(50, 5).Value = 0.002181059313589307
$ws.Cells.Item(51, 2).Value = 40.59999999999999
$ws.Cells.Item(51, 3).Value = 40.69865036010742
$ws.Cells.Item(51, 4).Value = 0.09865036010742756
$ws.Cells.Item(51, 5).Value = 0.009731893549325134
$ws.Cells.Item(52, 3).Value = 0.04815361022947684
$ws.Cells.Item(52, 5).Value = 2.654464649066489
$ws.Cells.Item(53, 5).Value = 0.05308929298132978
